# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
# Swap the "Periodo Mora" (col E) and "Valor Mora" (col F) values between
# the two worker rows (16 and 17) on the account-statement sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E16").Value = "1812"
$ws.Range("F16").Value = 32851

$ws.Range("E17").Value = "1901"
$ws.Range("F17").Value = 36502
